$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LiveData")

$ws.Range("C3").Value = 158806
$ws.Range("C4").Value = 149862
$ws.Range("C7").Value = 5.63
$ws.Range("C8").Value = 64.04000000000001
